$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (old Resolving-Mac sending-cluster rows removed entirely)
$ws.Range("A14:T17").EntireRow.Delete()

# Build new values for data rows 2-13 (updated TPM-based stats)
$arr = New-Object 'object[,]' 12,20
$arr[0,0] = "ECs"
$arr[0,1] = "Reln"
$arr[0,2] = "Lrp8"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 0.1471086666666667
$arr[0,7] = 0.441326
$arr[0,8] = 0.03503939655440032
$arr[0,9] = 0.03503939655440032
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 0.4394373333333333
$arr[0,13] = 1.318312
$arr[0,14] = 0.2944741752765458
$arr[0,15] = 0.2944741752765458
$arr[0,16] = 0.06464504019022221
$arr[0,17] = 0.5818053617119999
$arr[0,18] = 0.01031819740254488
$arr[0,19] = 0.01031819740254488
$arr[1,0] = "ECs"
$arr[1,1] = "Reln"
$arr[1,2] = "Lrp8"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 0.1471086666666667
$arr[1,7] = 0.441326
$arr[1,8] = 0.03503939655440032
$arr[1,9] = 0.03503939655440032
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 0.6503416666666667
$arr[1,13] = 1.951025
$arr[1,14] = 0.4358046333636673
$arr[1,15] = 0.4358046333636673
$arr[1,16] = 0.09567089546111111
$arr[1,17] = 0.8610380591500001
$arr[1,18] = 0.01527033136867458
$arr[1,19] = 0.01527033136867458
$arr[2,0] = "ECs"
$arr[2,1] = "Reln"
$arr[2,2] = "Lrp8"
$arr[2,3] = "MuSCs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.1471086666666667
$arr[2,7] = 0.441326
$arr[2,8] = 0.03503939655440032
$arr[2,9] = 0.03503939655440032
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 0.1607546666666667
$arr[2,13] = 0.482264
$arr[2,14] = 0.1077243426939663
$arr[2,15] = 0.1077243426939663
$arr[2,16] = 0.02364840467377778
$arr[2,17] = 0.212835642064
$arr[2,18] = 0.003774595962216003
$arr[2,19] = 0.003774595962216003
$arr[3,0] = "ECs"
$arr[3,1] = "Reln"
$arr[3,2] = "Lrp8"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.1471086666666667
$arr[3,7] = 0.441326
$arr[3,8] = 0.03503939655440032
$arr[3,9] = 0.03503939655440032
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 0.2417443333333333
$arr[3,13] = 0.725233
$arr[3,14] = 0.1619968486658205
$arr[3,15] = 0.1619968486658205
$arr[3,16] = 0.03556268655088889
$arr[3,17] = 0.320064178958
$arr[3,18] = 0.005676271820964862
$arr[3,19] = 0.005676271820964862
$arr[4,0] = "FAPs"
$arr[4,1] = "Reln"
$arr[4,2] = "Lrp8"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 1.029733
$arr[4,7] = 3.089199
$arr[4,8] = 0.2452691860358486
$arr[4,9] = 0.2452691860358485
$arr[4,10] = 2
$arr[4,11] = 0.6666666666666666
$arr[4,12] = 0.4394373333333333
$arr[4,13] = 1.318312
$arr[4,14] = 0.2944741752765458
$arr[4,15] = 0.2944741752765458
$arr[4,16] = 0.4525031235653333
$arr[4,17] = 4.072528112087999
$arr[4,18] = 0.07222544127865621
$arr[4,19] = 0.0722254412786562
$arr[5,0] = "FAPs"
$arr[5,1] = "Reln"
$arr[5,2] = "Lrp8"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 1.029733
$arr[5,7] = 3.089199
$arr[5,8] = 0.2452691860358486
$arr[5,9] = 0.2452691860358485
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 0.6503416666666667
$arr[5,13] = 1.951025
$arr[5,14] = 0.4358046333636673
$arr[5,15] = 0.4358046333636673
$arr[5,16] = 0.6696782754416667
$arr[5,17] = 6.027104478975001
$arr[5,18] = 0.1068894476957581
$arr[5,19] = 0.1068894476957581
$arr[6,0] = "FAPs"
$arr[6,1] = "Reln"
$arr[6,2] = "Lrp8"
$arr[6,3] = "MuSCs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 1.029733
$arr[6,7] = 3.089199
$arr[6,8] = 0.2452691860358486
$arr[6,9] = 0.2452691860358485
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 0.1607546666666667
$arr[6,13] = 0.482264
$arr[6,14] = 0.1077243426939663
$arr[6,15] = 0.1077243426939663
$arr[6,16] = 0.1655343851706667
$arr[6,17] = 1.489809466536
$arr[6,18] = 0.02642146184879594
$arr[6,19] = 0.02642146184879593
$arr[7,0] = "FAPs"
$arr[7,1] = "Reln"
$arr[7,2] = "Lrp8"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 1.029733
$arr[7,7] = 3.089199
$arr[7,8] = 0.2452691860358486
$arr[7,9] = 0.2452691860358485
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 0.2417443333333333
$arr[7,13] = 0.725233
$arr[7,14] = 0.1619968486658205
$arr[7,15] = 0.1619968486658205
$arr[7,16] = 0.2489321175963334
$arr[7,17] = 2.240389058367
$arr[7,18] = 0.03973283521263835
$arr[7,19] = 0.03973283521263834
$arr[8,0] = "MuSCs"
$arr[8,1] = "Reln"
$arr[8,2] = "Lrp8"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 3.021537333333333
$arr[8,7] = 9.064612
$arr[8,8] = 0.7196914174097511
$arr[8,9] = 0.7196914174097511
$arr[8,10] = 2
$arr[8,11] = 0.6666666666666666
$arr[8,12] = 0.4394373333333333
$arr[8,13] = 1.318312
$arr[8,14] = 0.2944741752765458
$arr[8,15] = 0.2944741752765458
$arr[8,16] = 1.327776308327111
$arr[8,17] = 11.949986774944
$arr[8,18] = 0.2119305365953447
$arr[8,19] = 0.2119305365953447
$arr[9,0] = "MuSCs"
$arr[9,1] = "Reln"
$arr[9,2] = "Lrp8"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 3.021537333333333
$arr[9,7] = 9.064612
$arr[9,8] = 0.7196914174097511
$arr[9,9] = 0.7196914174097511
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 0.6503416666666667
$arr[9,13] = 1.951025
$arr[9,14] = 0.4358046333636673
$arr[9,15] = 0.4358046333636673
$arr[9,16] = 1.965031625255556
$arr[9,17] = 17.6852846273
$arr[9,18] = 0.3136448542992346
$arr[9,19] = 0.3136448542992346
$arr[10,0] = "MuSCs"
$arr[10,1] = "Reln"
$arr[10,2] = "Lrp8"
$arr[10,3] = "MuSCs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 3.021537333333333
$arr[10,7] = 9.064612
$arr[10,8] = 0.7196914174097511
$arr[10,9] = 0.7196914174097511
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 0.1607546666666667
$arr[10,13] = 0.482264
$arr[10,14] = 0.1077243426939663
$arr[10,15] = 0.1077243426939663
$arr[10,16] = 0.485726226840889
$arr[10,17] = 4.371536041568
$arr[10,18] = 0.07752828488295439
$arr[10,19] = 0.07752828488295438
$arr[11,0] = "MuSCs"
$arr[11,1] = "Reln"
$arr[11,2] = "Lrp8"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 3.021537333333333
$arr[11,7] = 9.064612
$arr[11,8] = 0.7196914174097511
$arr[11,9] = 0.7196914174097511
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 0.2417443333333333
$arr[11,13] = 0.725233
$arr[11,14] = 0.1619968486658205
$arr[11,15] = 0.1619968486658205
$arr[11,16] = 0.7304395282884445
$arr[11,17] = 6.573955754596001
$arr[11,18] = 0.1165877416322173
$arr[11,19] = 0.1165877416322173

$ws.Range("A2:T13").Value = $arr
